# Refresh the cryptos list (prices + 1h volume deltas) for this run.
# D-column cells whose new value parses as a plain number are written via
# .Formula with a leading apostrophe (quote-prefix) so Excel keeps them as
# text, matching the original sheet where every price is a text string
# (several look like "60.235.98" which can't round-trip through a double).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.235.98"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.601.88"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = "'575.82"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").Formula = "'142.82"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Formula = "'0.999"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Formula = "'0.599"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("D9").Value = "2.609.32"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Formula = "'6.56"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Formula = "'0.105"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "3.061.33"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Formula = "'24.28"
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "60.232.55"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "2.604.25"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Formula = "'11.35"
$ws.Range("E19").Value = "  +6.34%  "
$ws.Range("D20").Formula = "'4.62"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Formula = "'346.40"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Formula = "'6.89"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Formula = "'63.10"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Formula = "'8.00"
$ws.Range("E28").Value = "  +4.30%  "
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D31").Formula = "'6.38"
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Formula = "'166.34"
$ws.Range("E33").Value = "  +4.80%  "
$ws.Range("D34").Formula = "'19.42"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Formula = "'1.30"
$ws.Range("E35").Value = "  +9.03%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Formula = "'4.27"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Formula = "'0.983"
$ws.Range("E37").Value = "  +6.08%  "
$ws.Range("E38").Value = "  +6.28%  "
$ws.Range("D39").Formula = "'38.06"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Formula = "'311.84"
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Formula = "'0.837"
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").Formula = "'135.66"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("D44").Formula = "'0.0995"
$ws.Range("E44").Value = "  +1.51%  "
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Formula = "'19.83"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").Formula = "'4.98"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Formula = "'19.95"
$ws.Range("E51").Value = "  +4.63%  "
